# "changes to git basics"
#
# The deck currently has a single slide ("Git Commands" with a screenshot).
# This adds a second slide ("Git bash (command line )") that uses the same
# Title+Content layout as slide 1 (slideLayout2.xml / "Title and Content"),
# with a bullet list of git-bash commands in the content placeholder.

$p = $ppt.ActivePresentation

# Insert the new slide as slide #2, using the Title+Content autolayout
# (same layout slide 1 already uses).
$s2 = $p.Slides.Add(2, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText)

# --- Title placeholder ---------------------------------------------------
$title = $s2.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Git bash (command line )"
$title.LanguageID = "en-GB"

# --- Content placeholder ---------------------------------------------------
# Build the bullet list one paragraph at a time so each paragraph's
# run gets its own language tag (matching the source deck's authoring).
$lines = @(
    "Git pull",
    "Git commit –m “comment”",
    "Git push –u origin master",
    "Git add .",
    "Git checkout <new Branch>",
    "Git merge <branch>",
    "Touch <new file> ( to create files) in local"
)

$body = $s2.Shapes.Item(2).TextFrame.TextRange
$body.Text = $lines[0]
$body.LanguageID = "en-GB"

for ($i = 1; $i -lt $lines.Count; $i++) {
    $body.InsertAfter("`r" + $lines[$i]) | Out-Null
    $para = $body.Paragraphs($i + 1, 1)
    $para.LanguageID = "en-GB"
}
